$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric-looking strings (e.g. "0.510", "19.90")
# must be forced to stay as TEXT (matching the source workbook, which stores
# every data cell as an inline/shared string) instead of being auto-coerced
# into a floating point number by Excel's usual "smart" cell-value parsing.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '26.884.36'
$ws.Range('E2').Value = '  +0.28%  '
$ws.Range('D3').Value = '1.638.62'
$ws.Range('E3').Value = '  -0.16%  '
$ws.Range('E4').Value = '  -0.57%  '
Set-TextValue $ws.Range('D5') '216.97'
$ws.Range('E5').Value = '  -0.67%  '
Set-TextValue $ws.Range('D6') '0.510'
$ws.Range('E6').Value = '  +1.93%  '
$ws.Range('E7').Value = '  -0.56%  '
$ws.Range('E8').Value = '  +1.70%  '
Set-TextValue $ws.Range('D9') '0.0625'
$ws.Range('E9').Value = '  +0.54%  '
Set-TextValue $ws.Range('D10') '19.90'
$ws.Range('E10').Value = '  +3.80%  '
Set-TextValue $ws.Range('D11') '0.0849'
$ws.Range('E11').Value = '  +0.32%  '
$ws.Range('D12').Value = '1.868.58'
$ws.Range('E12').Value = '  -0.13%  '
$ws.Range('D13').Value = '1.632.56'
$ws.Range('E13').Value = '  -0.54%  '
$ws.Range('E14').Value = '  -0.70%  '
$ws.Range('E15').Value = '  +1.08%  '
Set-TextValue $ws.Range('D16') '67.20'
$ws.Range('E16').Value = '  +3.18%  '
$ws.Range('D17').Value = '26.884.21'
$ws.Range('E17').Value = '  +0.22%  '
$ws.Range('D18').Value = '0.0₃0731'
$ws.Range('E18').Value = '  -0.24%  '
Set-TextValue $ws.Range('D19') '219.39'
$ws.Range('E19').Value = '  +1.70%  '
Set-TextValue $ws.Range('D21') '6.86'
$ws.Range('E21').Value = '  +3.75%  '
Set-TextValue $ws.Range('D22') '4.40'
$ws.Range('E22').Value = '  +0.80%  '
$ws.Range('E23').Value = '  +3.97%  '
Set-TextValue $ws.Range('D24') '9.16'
$ws.Range('E24').Value = '  +0.02%  '
Set-TextValue $ws.Range('D25') '146.87'
$ws.Range('E25').Value = '  -0.45%  '
$ws.Range('E26').Value = '  -0.66%  '
$ws.Range('E27').Value = '  +3.36%  '
$ws.Range('E28').Value = '  +0.71%  '
$ws.Range('E29').Value = '  +0.52%  '
Set-TextValue $ws.Range('D30') '0.0504'
$ws.Range('E30').Value = '  -1.09%  '
$ws.Range('E31').Value = '  -0.76%  '
$ws.Range('E32').Value = '  -1.35%  '
Set-TextValue $ws.Range('D33') '3.00'
$ws.Range('E33').Value = '  +0.80%  '
$ws.Range('D35').Value = '1.259.83'
$ws.Range('E35').Value = '  -0.44%  '
Set-TextValue $ws.Range('D36') '2.43'
$ws.Range('E36').Value = '  -0.34%  '
$ws.Range('E37').Value = '  +2.39%  '
Set-TextValue $ws.Range('D38') '0.535'
$ws.Range('E38').Value = '  +0.85%  '
$ws.Range('E39').Value = '  +2.18%  '
$ws.Range('E40').Value = '  -0.59%  '
$ws.Range('E41').Value = '  +0.84%  '
$ws.Range('E42').Value = '  +1.22%  '
$ws.Range('D43').Value = '1.778.75'
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range('D44') '2.10'
$ws.Range('E44').Value = '  -1.57%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D45') '61.75'
$ws.Range('E45').Value = '  +1.48%  '
Set-TextValue $ws.Range('D46') '91.65'
$ws.Range('E46').Value = '  -1.12%  '
$ws.Range('E47').Value = '  -1.11%  '
$ws.Range('E48').Value = '  +3.65%  '
$ws.Range('E49').Value = '  -0.20%  '
Set-TextValue $ws.Range('D50') '7.66'
$ws.Range('E50').Value = '  +1.41%  '
Set-TextValue $ws.Range('D51') '0.0963'
$ws.Range('E51').Value = '  -0.27%  '
